# Fixed a bug in mask
# Reorders the data rows (A2:F25) of the active sheet to the corrected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901,  16, 15, 45, 60, 60),
    @(701,  3,  90, 45, 97, 15),
    @(601,  9,  60, 67, 60, 42),
    @(1202, 2,  10, 10, 10, 10),
    @(1203, 3,  15, 15, 15, 15),
    @(101,  9,  30, 15, 60, 15),
    @(902,  1,  0,  0,  0,  0),
    @(501,  9,  52, 30, 75, 45),
    @(401,  9,  48, 67, 75, 45),
    @(201,  9,  30, 15, 45, 30),
    @(1201, 2,  10, 10, 10, 10),
    @(1001, 18, 30, 75, 60, 72),
    @(301,  6,  45, 30, 60, 45),
    @(801,  3,  67, 65, 52, 45),
    @(1,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0,  15, 30, 30, 0),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(802,  0,  4,  5,  4,  0),
    @(602,  0,  0,  4,  0,  9),
    @(402,  0,  0,  4,  0,  0),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
